# Update the table style used by every data table in the deck (the "Update
# Data Sources from LFX" refresh restyles all data tables to the new table
# style id).
$OldStyleId = "{C13E3245-A485-425E-8DFC-A9FE81DB0A83}"
$NewStyleId = "{7E67995A-638B-42D5-9D7F-B0D75F3B4303}"

$p = $ppt.ActivePresentation

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shp = $slide.Shapes.Item($j)
        if ($shp.HasTable) {
            $tbl = $shp.Table
            if ($tbl.Style -eq $OldStyleId) {
                $tbl.ApplyStyle($NewStyleId)
            }
        }
    }
}
